# Auto-generated edit script: applies the scraped-data refresh for Linea 141 (12/01/2026)
# Updates "Ultima actualizacion" / "Total filas" headers and refreshes/appends Hora_Scrap rows
# across the three worksheets (LP1912, LP1912-215, 6203-6173).

$wb = $excel.ActiveWorkbook

# ---- Worksheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")

$ws.Cells.Item(2, 1).Value = "Última actualización: 16:14:52"
$ws.Cells.Item(3, 1).Value = "Total filas: 366"

# Hora_Scrap / Hora_Llegada / Linea / Minutos / Parada (columns A:E)
$rows = @(
    @(73, "06:52:23", "08:23", "215B_EL PATO", 91, "LP1912"),
    @(74, "07:46:15", "08:23", "16_P MOR-SANTA ANA", 37, "LP1912"),
    @(118, "08:39:56", "09:34", "16_SANTA ANA", 55, "LP1912"),
    @(119, "08:39:56", "09:34", "23_HERNANDEZ", 55, "LP1912"),
    @(142, "08:50:00", "10:27", "215A_EL PATO", 97, "LP1912"),
    @(143, "09:38:04", "10:27", "26_HERNANDEZ", 49, "LP1912"),
    @(152, "10:28:12", "10:57", "23_HERNANDEZ", 29, "LP1912"),
    @(153, "10:57:58", "10:57", "17_ROMERO", 0, "LP1912"),
    @(154, "10:28:12", "10:57", "27_EL RETIRO", 29, "LP1912"),
    @(165, "09:38:04", "11:12", "23_HERNANDEZ", 94, "LP1912"),
    @(166, "10:28:12", "11:12", "15_ABASTO", 44, "LP1912"),
    @(176, "10:57:58", "11:26", "23_HERNANDEZ", 29, "LP1912"),
    @(177, "10:57:58", "11:26", "225_C ROCA-H SUR", 29, "LP1912"),
    @(199, "10:57:58", "12:06", "16_P MOR-SANTA ANA", 69, "LP1912"),
    @(200, "10:57:58", "12:06", "14_ABASTO", 69, "LP1912"),
    @(201, "10:28:12", "12:06", "84_COLONIA URQUIZA-ESC 49", 98, "LP1912"),
    @(208, "10:57:58", "12:14", "10_OLMOS", 77, "LP1912"),
    @(209, "10:28:12", "12:14", "17_ROMERO", 106, "LP1912"),
    @(222, "10:57:58", "12:36", "27_EL RETIRO", 99, "LP1912"),
    @(223, "12:16:51", "12:36", "16_SANTA ANA", 20, "LP1912"),
    @(280, "12:44:21", "14:20", "215C_EL PATO", 96, "LP1912"),
    @(281, "14:16:51", "14:20", "26_HERNANDEZ", 4, "LP1912"),
    @(303, "14:40:41", "15:13", "10_OLMOS", 33, "LP1912"),
    @(304, "14:16:51", "15:13", "11_ETCHEVERRY", 57, "LP1912"),
    @(341, "16:14:52", "16:22", "14_ABASTO", 8, "LP1912"),
    @(342, "16:14:52", "16:26", "16_SANTA ANA", 12, "LP1912"),
    @(343, "16:14:52", "16:28", "10_OLMOS", 14, "LP1912"),
    @(344, "15:51:40", "16:29", "10_OLMOS", 38, "LP1912"),
    @(345, "14:40:41", "16:30", "15_ABASTO", 110, "LP1912"),
    @(346, "15:51:40", "16:34", "23_HERNANDEZ", 43, "LP1912"),
    @(347, "16:14:52", "16:34", "16_SANTA ANA", 20, "LP1912"),
    @(348, "16:14:52", "16:35", "23_HERNANDEZ", 21, "LP1912"),
    @(349, "15:19:52", "16:36", "11_ETCHEVERRY", 77, "LP1912"),
    @(350, "15:19:52", "16:39", "17_ROMERO", 80, "LP1912"),
    @(351, "14:53:55", "16:42", "16_P MOR-SANTA ANA", 109, "LP1912"),
    @(352, "14:53:55", "16:42", "225_GOMEZ", 109, "LP1912"),
    @(353, "15:51:40", "16:43", "225_GOMEZ", 52, "LP1912"),
    @(354, "14:53:55", "16:48", "15_ABASTO", 115, "LP1912"),
    @(355, "15:51:40", "16:50", "14_ABASTO", 59, "LP1912"),
    @(356, "15:19:52", "16:56", "17_179 Y 38", 97, "LP1912"),
    @(357, "16:14:52", "16:56", "10_OLMOS", 42, "LP1912"),
    @(358, "16:14:52", "17:04", "23_HERNANDEZ", 50, "LP1912"),
    @(359, "15:51:40", "17:04", "11_ETCHEVERRY", 73, "LP1912"),
    @(360, "15:19:52", "17:04", "215A_EL PATO", 105, "LP1912"),
    @(361, "16:14:52", "17:20", "26_HERNANDEZ", 66, "LP1912"),
    @(362, "15:51:40", "17:21", "26_HERNANDEZ", 90, "LP1912"),
    @(363, "15:51:40", "17:24", "84_COLONIA URQUIZA-ESC 49", 93, "LP1912"),
    @(364, "15:51:40", "17:28", "14_ABASTO", 97, "LP1912"),
    @(365, "16:14:52", "17:35", "27_EL RETIRO", 81, "LP1912"),
    @(366, "15:51:40", "17:36", "27_EL RETIRO", 105, "LP1912"),
    @(367, "15:51:40", "17:38", "17_ROMERO", 107, "LP1912"),
    @(368, "15:51:40", "17:40", "215B_EL PATO", 109, "LP1912"),
    @(369, "15:51:40", "17:50", "16_P MOR-167 Y 521", 119, "LP1912"),
    @(370, "16:14:52", "17:52", "81_EL PELIGRO", 98, "LP1912"),
    @(371, "16:14:52", "18:04", "17_ROMERO", 110, "LP1912")
)
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}


# ---- Worksheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")

$ws.Cells.Item(2, 1).Value = "Última actualización: 16:14:52"


# ---- Worksheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")

$ws.Cells.Item(2, 1).Value = "Última actualización: 16:14:52"
$ws.Cells.Item(3, 1).Value = "Total filas: 50"

# Hora_Scrap / Hora_Llegada / Linea / Minutos / Parada (columns A:E)
$rows = @(
    @(50, "16:14:52", "16:14", "215C_LA PLATA", 0, "L6203"),
    @(51, "15:19:52", "16:52", "215B_LP-P MOR-40 Y 115", 93, "L6173"),
    @(52, "15:19:52", "17:14", "215A_LA PLATA", 115, "L6173"),
    @(53, "15:51:40", "17:15", "215A_LA PLATA", 84, "L6173"),
    @(54, "16:14:52", "17:17", "215A_LA PLATA", 63, "L6173"),
    @(55, "16:14:52", "18:03", "215C_LA PLATA", 109, "L6203")
)
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

